$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete entire row 5 (Resolving-Mac row), which also removes the now-unused
# "Resolving-Mac" shared string and shifts later indices down.
$ws.Rows.Item(5).Delete() | Out-Null

# Update numeric values for row 2 (ECs -> Qrfp/Qrfpr)
$ws.Range("G2").Value = 0.2581393333333333
$ws.Range("H2").Value = 0.7744180000000001
$ws.Range("I2").Value = 0.174859595118225
$ws.Range("J2").Value = 0.1748595951182251
$ws.Range("Q2").Value = 0.03132667088955556
$ws.Range("R2").Value = 0.281940038006
$ws.Range("S2").Value = 0.174859595118225
$ws.Range("T2").Value = 0.1748595951182251

# Update numeric values for row 3 (FAPs -> Qrfp/Qrfpr)
$ws.Range("I3").Value = 0.5237451507733812
$ws.Range("J3").Value = 0.5237451507733814
$ws.Range("Q3").Value = 0.0938306642948889
$ws.Range("R3").Value = 0.844475978654
$ws.Range("S3").Value = 0.5237451507733812
$ws.Range("T3").Value = 0.5237451507733814

# Update numeric values for row 4 (MuSCs -> Qrfp/Qrfpr)
$ws.Range("G4").Value = 0.4449396666666667
$ws.Range("H4").Value = 1.334819
$ws.Range("I4").Value = 0.3013952541083937
$ws.Range("J4").Value = 0.3013952541083937
$ws.Range("Q4").Value = 0.05399594987477777
$ws.Range("R4").Value = 0.485963548873
$ws.Range("S4").Value = 0.3013952541083937
$ws.Range("T4").Value = 0.3013952541083937

$wb.Save()
